# Update PLC data 2025-10-13 13:40:57
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = 153970
$ws.Range("C5").Value = 8650
$ws.Range("C6").Value = 515
$ws.Range("C7").Value = 5.62
